$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5125
$ws.Range("K3").Value = 5267
$ws.Range("K4").Value = 1099
$ws.Range("K5").Value = 377
$ws.Range("K6").Value = 5907
$ws.Range("K7").Value = 17775

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K7").Value = 532
$ws.Range("K8").Value = 1204
$ws.Range("K10").Value = 98
$ws.Range("K15").Value = 179
$ws.Range("K17").Value = 34
$ws.Range("K19").Value = 527
$ws.Range("K20").Value = 407
$ws.Range("K21").Value = 55
$ws.Range("K22").Value = 45
$ws.Range("K23").Value = 182
$ws.Range("K25").Value = 85
$ws.Range("K29").Value = 955
$ws.Range("K33").Value = 754
$ws.Range("K37").Value = 600
$ws.Range("K42").Value = 654
$ws.Range("K43").Value = 154
$ws.Range("K45").Value = 20
$ws.Range("K48").Value = 221
$ws.Range("K51").Value = 226
$ws.Range("K52").Value = 465
$ws.Range("K54").Value = 348
$ws.Range("K55").Value = 200
$ws.Range("K57").Value = 65
$ws.Range("K60").Value = 111
$ws.Range("K63").Value = 55
$ws.Range("K64").Value = 111
$ws.Range("K65").Value = 404
$ws.Range("K67").Value = 678
$ws.Range("K69").Value = 39
$ws.Range("K76").Value = 247
$ws.Range("K78").Value = 203
$ws.Range("K79").Value = 436
$ws.Range("K80").Value = 64
$ws.Range("K83").Value = 395
$ws.Range("K85").Value = 830
$ws.Range("K87").Value = 30
$ws.Range("K88").Value = 197
$ws.Range("K89").Value = 259
$ws.Range("K90").Value = 159
$ws.Range("K95").Value = 306
$ws.Range("K96").Value = 191
$ws.Range("K99").Value = 302
$ws.Range("K101").Value = 17775

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K6").Value = 82
$ws.Range("K7").Value = 191

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K5").Value = 21
$ws.Range("K6").Value = 138
$ws.Range("K7").Value = 532

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K2").Value = 72
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 259

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 280
$ws.Range("K3").Value = 279
$ws.Range("K6").Value = 197
$ws.Range("K7").Value = 830

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K3").Value = 129
$ws.Range("K6").Value = 172
$ws.Range("K7").Value = 465

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 331
$ws.Range("K3").Value = 361
$ws.Range("K4").Value = 69
$ws.Range("K6").Value = 407
$ws.Range("K7").Value = 1204

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 134
$ws.Range("K7").Value = 395

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K2").Value = 203
$ws.Range("K3").Value = 281
$ws.Range("K7").Value = 754

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K2").Value = 102
$ws.Range("K3").Value = 109
$ws.Range("K6").Value = 70
$ws.Range("K7").Value = 306

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 197
$ws.Range("K7").Value = 600

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K6").Value = 153
$ws.Range("K7").Value = 404

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 125
$ws.Range("K7").Value = 302

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 197
$ws.Range("K3").Value = 238
$ws.Range("K7").Value = 678

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K2").Value = 19
$ws.Range("K3").Value = 18

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K2").Value = 54
$ws.Range("K6").Value = 184
$ws.Range("K7").Value = 348

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K3").Value = 340
$ws.Range("K6").Value = 266
$ws.Range("K7").Value = 955

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 221

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K2").Value = 155
$ws.Range("K4").Value = 23
$ws.Range("K7").Value = 527

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 247

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K6").Value = 248
$ws.Range("K7").Value = 654

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 73
$ws.Range("K7").Value = 203

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 200

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K2").Value = 49
$ws.Range("K6").Value = 49
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K6").Value = 109
$ws.Range("K7").Value = 436

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 135
$ws.Range("K6").Value = 117
$ws.Range("K7").Value = 407

$ws = $wb.Worksheets.Item('Burnside')
$ws.Range("K4").Value = 3
$ws.Range("K7").Value = 34

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 28
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 197

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K3").Value = 60
$ws.Range("K7").Value = 226

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K3").Value = 33
$ws.Range("K7").Value = 111

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Jackson Park')
$ws.Range("K3").Value = 6
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 30
